$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 10: fill in previously-blank placeholder row (keeps its existing style) ---
$ws.Cells.Item(10, 1).Value = 9
$ws.Cells.Item(10, 2).Value = 141
$ws.Cells.Item(10, 3).Value = "给定一个链表，判断链表中是否有环"
$ws.Cells.Item(10, 4).Value = "快慢指针法`n1 快慢指针初始指向头节点`n2 快指针&&快指针的next不为空`n2 快指针一次走两步，慢指针一次走一步`n3 判断是否相等，如果相等代表有环，快慢指针相遇，返回结果，退出程序；否则没有环`n4 循环结束，没有找到相等的节点，就代表没有相交点，返回false"
$ws.Cells.Item(10, 5).Value = "快慢指针"
$ws.Cells.Item(10, 5).VerticalAlignment = -4160
$ws.Cells.Item(10, 6).Value = "1）没有环：快指针走到链表尾部，就退出循环，慢指针访问了链表一半的数据，快指针访问了所有的数据，O(n/2) + O(n)，时间复杂度是O(n)`n2）有环：慢指针进入环之前，走过的长度是M，快指针在环中迭代的元素个数是M，即非环长度是M。`n慢指针进入环之后，快慢指针经过多长时间会相遇（经过多少次迭代）：快慢指针的距离(最大是环的长度N)/快慢指针差即1，遍历迭代次数O(M + N), 即最大是O(链表长度)"
$ws.Cells.Item(10, 7).Value = "O(1)"
$ws.Rows.Item(10).RowHeight = 260

# --- Row 11: new row (picks up default column style) ---
$ws.Cells.Item(11, 1).Value = 10
$ws.Cells.Item(11, 2).Value = 142
$ws.Cells.Item(11, 3).Value = "给定一个链表，返回链表开始入环的第一个节点。 如果链表无环，则返回 null"
$ws.Cells.Item(11, 4).Value = "1 快慢指针同时指向链表头节点`n2 快指针以及后续指针不为空`n3 快指针走两步，慢指针走一步`n3 快慢指针对应的节点是否相等，如果相等就说明有环，进一步验证环的入口（有环就一定有入口）`n，快节点从头开始走，慢节点继续前进，步长均为1，当slow==fast，说明入口到了`n4 如果没有找到或者快节点为空，就说明没有环，返回null"
$ws.Cells.Item(11, 5).Value = "快慢指针"
$ws.Cells.Item(11, 7).Value = "O(1)"
$ws.Rows.Item(11).RowHeight = 180

# --- Row 12: new row (picks up default column style); note D before C to match shared-string order ---
$ws.Cells.Item(12, 1).Value = 11
$ws.Cells.Item(12, 2).Value = 143
$ws.Cells.Item(12, 4).Value = "1 计算链表的中间节点：`n     链表长度是奇数，slow是链表中间节点`n      链表长度是偶数，slow是链表中间位置偏右侧的节点`n2 链表的后半部分反转，得到反转后的链表头部`n3 链表前半部分与反转后的链表穿插串联节点"
$ws.Cells.Item(12, 3).Value = "//给定一个单链表 L：L0→L1→…→Ln-1→Ln ， `n//将其重新排列后变为： L0→Ln→L1→Ln-1→L2→Ln-2→… "
$ws.Cells.Item(12, 5).Value = "快慢指针`n链表反转`n两个链表拼接"
$ws.Cells.Item(12, 6).Value = "O(M+N),M,N是两个链表的元素个数"
$ws.Cells.Item(12, 7).Value = "O(1)"
$ws.Rows.Item(12).RowHeight = 120

# --- View state: final selection lands on F13 ---
$ws.Range("F13").Select()
